$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (e.g. "17.70") stay as text, matching source data
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D18", "D20", "D27", "D28", "D33", "D35", "D37", "D42", "D44", "D46", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.165.25"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "2.269.96"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "498.48"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "128.18"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "0.0952"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +3.49%  "
$ws.Range("D12").Value = "4.74"
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").Value = "2.670.49"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "22.65"
$ws.Range("E14").Value = "  +4.07%  "
$ws.Range("D15").Value = "54.144.51"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "2.267.82"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "10.21"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "302.79"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").Value = "170.76"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "1.60"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "0.0₃0692"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D33").Value = "17.70"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "0.942"
$ws.Range("E35").Value = "  +9.78%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "3.70"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "124.60"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("D44").Value = "0.0889"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "238.55"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "10.76"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("E51").Value = "  -0.46%  "
